# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values. Update the rows whose K (strikeout) counts
# changed when the save data was regenerated.
$gValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 3
    12 = 2
    13 = 2
    16 = 4
    17 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
